$d = $word.ActiveDocument

# --- Add the three new character styles ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Waktu Kampanye 2022 ..." run (4 occurrences) ---

$campaignText = "Waktu Kampanye 2022 untuk Rasi bintang Pegasus: 8-17 Oktober, 7-16 November,"
$rng = $d.Content
$rng.Start = 0
$found = $rng.Find.Execute($campaignText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
    $found = $rng.Find.Execute($campaignText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the participation paragraph ---

$paragraphText = "Anda sedang berpartisipasi dalam kampanye global pengamatan dan pencatatan penampakan bintang paling redup untuk pengukuran tingkat polusi cahaya di suatu lokasi. Melalui pengamatan dan identifikasi  Rasi bintang Pegasus di langit malam dan membandingkannya dengan peta bintang, masyarakat di seluruh dunia dapat mengetahui dan mempelajari seberapa besar kontribusi cahaya di lingkungannya terhadap polusi cahaya. Kontribusi data anda pada basis data online akan membantu mendokumentasikan langit malam yang tampak di berbagai lokasi."
$rngPara = $d.Content
$foundPara = $rngPara.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPara) {
    $rngPara.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the map-credit paragraph ---

$linksText = "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rngLinks = $d.Content
$foundLinks = $rngLinks.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundLinks) {
    $rngLinks.Style = "GaNLinks"
}

Write-Output "Styles added and applied."
